$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 44-45, pushing the existing records (old rows 44-150) down to rows 46-152
$ws.Rows("44:45").Insert()

# Populate new row 44: new record (Especial)
$ws.Range("A44").Value = 1
$ws.Range("B44").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C44").Value = 'Arica y Parinacota'
$ws.Range("D44").Value = 44972
$ws.Range("E44").Value = 15
$ws.Range("F44").Value = 'Fruta'
$ws.Range("G44").Value = 100108
$ws.Range("H44").Value = 'Tropicales y subtropicales'
$ws.Range("I44").Value = 100108003
$ws.Range("J44").Value = 'Maracuyá'
$ws.Range("K44").Value = 'Sin especificar'
$ws.Range("L44").Value = 'Especial'
$ws.Range("M44").Value = 80
$ws.Range("N44").Value = 44000
$ws.Range("O44").Value = 45000
$ws.Range("P44").Value = 44625
$ws.Range("Q44").Value = '$/caja 20 kilos'
$ws.Range("R44").Value = 'Región de Arica y Parinacota'
$ws.Range("S44").Value = 2231
$ws.Range("T44").Value = 20

# Populate new row 45: new record (Primera)
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C45").Value = 'Arica y Parinacota'
$ws.Range("D45").Value = 44972
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = 'Fruta'
$ws.Range("G45").Value = 100108
$ws.Range("H45").Value = 'Tropicales y subtropicales'
$ws.Range("I45").Value = 100108003
$ws.Range("J45").Value = 'Maracuyá'
$ws.Range("K45").Value = 'Sin especificar'
$ws.Range("L45").Value = 'Primera'
$ws.Range("M45").Value = 90
$ws.Range("N45").Value = 39000
$ws.Range("O45").Value = 40000
$ws.Range("P45").Value = 39667
$ws.Range("Q45").Value = '$/caja 20 kilos'
$ws.Range("R45").Value = 'Región de Arica y Parinacota'
$ws.Range("S45").Value = 1983
$ws.Range("T45").Value = 20
